$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("E3").Value = 16

# Row 15
$ws.Range("E15").Value = 124

# Row 17
$ws.Range("E17").Value = 82

# Row 18
$ws.Range("E18").Value = 78

# Row 26
$ws.Range("E26").Value = 18

# Row 28
$ws.Range("E28").Value = 9

# Row 34
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = 5
$ws.Range("H34").Value = 5

# Row 38
$ws.Range("E38").Value = 48

# Row 42
$ws.Range("E42").Value = 26

# Row 46
$ws.Range("E46").Value = 19
$ws.Range("F46").Value = 5
$ws.Range("H46").Value = 5

# Row 47
$ws.Range("E47").Value = 44

# Row 49
$ws.Range("E49").Value = 50

# Row 50
$ws.Range("E50").Value = 13
$ws.Range("F50").Value = 2
$ws.Range("H50").Value = 2

# Row 55
$ws.Range("E55").Value = 5

# Row 61
$ws.Range("E61").Value = 20

# Row 62
$ws.Range("E62").Value = 28

# Row 63
$ws.Range("E63").Value = 18

# Row 69
$ws.Range("E69").Value = 12

# Row 70
$ws.Range("E70").Value = 25
$ws.Range("F70").Value = 12
$ws.Range("H70").Value = 12

# Row 74
$ws.Range("E74").Value = 14

# Row 77
$ws.Range("E77").Value = 32

# Row 78
$ws.Range("E78").Value = 23
$ws.Range("F78").Value = 9
$ws.Range("H78").Value = 9
